$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the two new columns
$ws.Range("N1").Value = "facultad"
$ws.Range("O1").Value = "departamento"

# Add the new data for row 3 (row 2 stays blank in these columns)
$ws.Range("N3").Value = "Facultad de Producción y diseño"
$ws.Range("O3").Value = "Departamento de producción"

# Update the visible selection (matches the saved view state in the diff)
$ws.Range("G19").Select()
